# Slide 6 ("Templafy") -> "Content Placeholder 2" -> 7th paragraph currently
# reads "Distributed team". The author placed the cursor right after the
# word "Distributed " and typed "development ", turning the bullet into
# "Distributed development team" (run split into "Distributed development "
# + "team").
$p    = $ppt.ActivePresentation
$s    = $p.Slides.Item(6)
$shp  = $s.Shapes.Item("Content Placeholder 2")
$para = $shp.TextFrame.TextRange.Paragraphs(7)

# The leading "Distributed " is the first 12 characters (incl. trailing
# space) of the paragraph; insert the new word right after it.
$lead = $para.Characters(1, 12)
[void]$lead.InsertAfter("development ")
